# Applies the "Added most instructions to interpreter" edit to the
# instruction-set workbook:
#   * Renames the DJZ mnemonic cell text (adds the missing comma)
#   * Documents two groups of instructions that still need to be added
#   * Fixes up the border box around the DEBUG Rn / DEBUG Message rows
#     (the DEBUG Rn row becomes its own boxed row, DEBUG Message becomes
#     the header of a new boxed area that spans rows 66-71)
#   * Restores the selection near the bottom of the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------

function Set-MediumEdge($range, [int]$edge) {
    $b = $range.Borders.Item($edge)
    $b.LineStyle = 1
    $b.Weight = -4138
    $b.Color = 0
}

$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10
$xlCenter = -4108

function Set-BoxEdges($range, [bool]$left, [bool]$top, [bool]$bottom, [bool]$right) {
    if ($left)   { Set-MediumEdge $range $xlEdgeLeft }
    if ($top)    { Set-MediumEdge $range $xlEdgeTop }
    if ($bottom) { Set-MediumEdge $range $xlEdgeBottom }
    if ($right)  { Set-MediumEdge $range $xlEdgeRight }
}

# ---------------------------------------------------------------------
# 1. Fix the DJZ mnemonic text (missing comma)
# ---------------------------------------------------------------------

$ws.Range("A52").Value = "DJZ Rn, label"

# ---------------------------------------------------------------------
# 2. Document the instructions still to be added
# ---------------------------------------------------------------------

$ws.Range("A55").Value = "Bitwise AND, OR, XOR, COMP to be added"
$ws.Range("A56").Value = "PUSH and POP to be added"

# ---------------------------------------------------------------------
# 3. Row 65 (DEBUG Rn) becomes its own boxed row
# ---------------------------------------------------------------------

Set-BoxEdges $ws.Range("A65") $true  $true $true $false
Set-BoxEdges $ws.Range("B65") $false $true $true $false
Set-BoxEdges $ws.Range("C65") $false $true $true $true
$ws.Rows.Item(65).RowHeight = 15

# ---------------------------------------------------------------------
# 4. Row 66 (DEBUG Message) becomes the top of a new box
# ---------------------------------------------------------------------

Set-BoxEdges $ws.Range("A66") $true  $true $false $false
Set-BoxEdges $ws.Range("B66") $false $true $false $false
Set-BoxEdges $ws.Range("C66") $false $true $false $true

# ---------------------------------------------------------------------
# 5. Rows 67-70 are the middle of the box: only left/right edges,
#    the "." filler rows stay centered (first one also vertically).
# ---------------------------------------------------------------------

Set-BoxEdges $ws.Range("A67") $true $false $false $false
Set-BoxEdges $ws.Range("C67") $false $false $false $true

Set-BoxEdges $ws.Range("A68") $true $false $false $false
Set-BoxEdges $ws.Range("C68") $false $false $false $true
$ws.Range("B68").HorizontalAlignment = $xlCenter
$ws.Range("B68").VerticalAlignment = $xlCenter

Set-BoxEdges $ws.Range("A69") $true $false $false $false
Set-BoxEdges $ws.Range("C69") $false $false $false $true
$ws.Range("B69").HorizontalAlignment = $xlCenter

Set-BoxEdges $ws.Range("A70") $true $false $false $false
Set-BoxEdges $ws.Range("C70") $false $false $false $true
$ws.Range("B70").HorizontalAlignment = $xlCenter

# ---------------------------------------------------------------------
# 6. Row 71 closes the box at the bottom
# ---------------------------------------------------------------------

Set-BoxEdges $ws.Range("A71") $true  $false $true $false
Set-BoxEdges $ws.Range("B71") $false $false $true $false
Set-BoxEdges $ws.Range("C71") $false $false $true $true
$ws.Rows.Item(71).RowHeight = 15

# ---------------------------------------------------------------------
# 7. Restore the selection near the bottom of the sheet
# ---------------------------------------------------------------------

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 46
$ws.Range("A56").Select()
